$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatted column block (D:K) to a uniform width of 9
# (matches stored width "9" once Excel's 5px padding is subtracted back out).
$ws.Range("D1:K1").ColumnWidth = 8.14

# Add the new "2022" column (K) by cloning the formatting of column J's
# header/data rows (3-5), then filling in the new figures.
$ws.Range("J3:J5").Copy()
$ws.Range("K3").PasteSpecial(-4122)

$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 0.11705180708279034

# Restore the selection to match the authored state.
$ws.Range("J12").Select()
